# Freelance Timesheet - add two new work-log entries (GCal error-handling work)
# and refresh the saved view/selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20 and 21 were blank placeholder rows; give them the same
# number/date/border formatting as the rows directly above them before
# filling in the new values.
$ws.Range("A18:E18").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A19:E19").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New timesheet entry: 2023-03-04, 12:12 - 13:12
$ws.Range("A20").Value = 44989
$ws.Range("B20").Value = 0.5083333333333333
$ws.Range("C20").Value = 0.54999999999999993
$ws.Range("D20").Value = "Added support for a calendar for every tutor"
$ws.Range("E20").Value = 1

# New timesheet entry: 2023-03-05, 03:20 - 04:20 (long description -> taller row)
$ws.Range("A21").Value = 44990
$ws.Range("B21").Value = 0.1388888888888889
$ws.Range("C21").Value = 0.18055555555555555
$ws.Range("D21").Value = "Improved Google Calendar multiplicity support. Error handling, etc"
$ws.Range("E21").Value = 1

$ws.Rows.Item(21).RowHeight = 28

# Totals (E32/E33) are formulas and recalc automatically.

# Refresh the saved scroll position / selection.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D23").Select()
